$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Trends Status" - update values
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Range("B2").Value = 2
$ws1.Range("C2").Value = 10
$ws1.Range("D2").Value = 7.4
$ws1.Range("E2").Value = 17.5

$ws1.Range("B3").Value = 5
$ws1.Range("C3").Value = 11
$ws1.Range("D3").Value = 18.5
$ws1.Range("E3").Value = 19.3

$ws1.Range("B4").Value = 3
$ws1.Range("C4").Value = 32
$ws1.Range("D4").Value = 11.1
$ws1.Range("E4").Value = 56.1

$ws1.Range("B5").Value = 7
$ws1.Range("C5").Value = 1
$ws1.Range("D5").Value = 25.9
$ws1.Range("E5").Value = 1.8

$ws1.Range("B6").Value = 10
$ws1.Range("C6").Value = 3
$ws1.Range("D6").Value = 37

$ws1.Range("B7").Value = 47
$ws1.Range("C7").Value = 171

$ws1.Range("B8").Value = 362
$ws1.Range("C8").Value = 208

# ---------------------------------------------------------------------------
# Sheet 3: "Priority Status" - update values
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Priority Status")
$ws3.Range("B2").Value = 103
$ws3.Range("B3").Value = 286
$ws3.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# Sheet 4: "Species qualification" - update values
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("A2").Value = "SoIB Assessment"
$ws4.Range("B2").Value = 436

$ws4.Range("B3").Value = 74
$ws4.Range("C3").Value = 27

$ws4.Range("B4").Value = 228
$ws4.Range("C4").Value = 57

# ---------------------------------------------------------------------------
# Sheet 5 "High Priority break-up":
#   - add a brand-new sheet ("Major update - High Priority ") that keeps the
#     OLD content of this sheet, placed at the end of the workbook
#   - rename the existing sheet to "Interannual update - High Pri" and
#     replace its contents with the new breakdown
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("High Priority break-up")

# Capture the old values before they get overwritten below. (Read via
# Value2 - the plain Value getter stringifies to a placeholder in this
# interpreter when captured into a variable.)
$oldHeaderA = $ws5.Range("A1").Value2
$oldHeaderB = $ws5.Range("B1").Value2
$oldHeaderC = $ws5.Range("C1").Value2
$oldHeaderD = $ws5.Range("D1").Value2
$oldHeaderE = $ws5.Range("E1").Value2

# Create the new sheet and name it while it is still easy to address (right
# after creation) - the engine's object references bind to sheet position,
# so rename/read-back-by-name BEFORE moving, and never use the handle again
# after a Move().
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "Major update - High Priority "

$wsNew.Range("A1").Value = $oldHeaderA
$wsNew.Range("B1").Value = $oldHeaderB
$wsNew.Range("C1").Value = $oldHeaderC
$wsNew.Range("D1").Value = $oldHeaderD
$wsNew.Range("E1").Value = $oldHeaderE
$wsNew.Range("A1:E1").Font.Bold = $true
$wsNew.Range("A1:E1").HorizontalAlignment = -4108

$wsNew.Range("A2").Value = "Trend New"
$wsNew.Range("B2").Value = 4
$wsNew.Range("C2").Value = 21.1
$wsNew.Range("D2").Value = 4
$wsNew.Range("E2").Value = 21.1

$wsNew.Range("A3").Value = "IUCN"
$wsNew.Range("B3").Value = 15
$wsNew.Range("C3").Value = 78.90000000000001
$wsNew.Range("D3").Value = 15
$wsNew.Range("E3").Value = 78.90000000000001

# Move the new sheet to the end of the workbook (after the current last
# sheet). Do this via the name lookup, and re-fetch by name afterward since
# the old handle no longer tracks the moved sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew.Move($null, $lastSheet)

# Now rename the original sheet and replace its contents with the new data.
$ws5 = $wb.Worksheets.Item("High Priority break-up")
$ws5.Name = "Interannual update - High Pri"

$ws5.Range("A2").Value = "Trend New"
$ws5.Range("B2").Value = 73
$ws5.Range("C2").Value = 70.90000000000001
$ws5.Range("D2").Value = 73
$ws5.Range("E2").Value = 81.09999999999999

$ws5.Range("A3").Value = "Trend Different"
$ws5.Range("B3").Value = 1
$ws5.Range("C3").Value = 1
$ws5.Range("D3").Value = $null
$ws5.Range("E3").Value = $null

$ws5.Range("A4").Value = "IUCN"
$ws5.Range("B4").Value = 29
$ws5.Range("C4").Value = 28.2
$ws5.Range("D4").Value = 17
$ws5.Range("E4").Value = 18.9
